$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7537016272544861
$ws.Range("B1").Value = 1.380809187889099
$ws.Range("C1").Value = 4.9274001121521
$ws.Range("D1").Value = 1.996675133705139
$ws.Range("E1").Value = 1.257204532623291
